$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.02
$ws.Range("H2").Value = 3.15
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 2.55
$ws.Range("K2").Value = 2.07
$ws.Range("L2").Value = 3.9
$ws.Range("O2").Value = 1.26
$ws.Range("P2").Value = 3.08
$ws.Range("Q2").Value = 1.91
$ws.Range("U2").Value = 1.79
$ws.Range("V2").Value = 1.98
$ws.Range("W2").Value = 6.2
$ws.Range("X2").Value = 8.25
$ws.Range("Y2").Value = 7.2
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 13
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 8.75
$ws.Range("AD2").Value = 5.4
$ws.Range("AE2").Value = 11.75
$ws.Range("AF2").Value = 50
$ws.Range("AG2").Value = 350
$ws.Range("AH2").Value = 8
$ws.Range("AI2").Value = 14.5
$ws.Range("AJ2").Value = 10
$ws.Range("AK2").Value = 37
$ws.Range("AL2").Value = 25
$ws.Range("AM2").Value = 30
$ws.Range("AN2").Value = 3.9
$ws.Range("AO2").Value = 10
$ws.Range("AP2").Value = 17.5
$ws.Range("AQ2").Value = 37
$ws.Range("AR2").Value = 65
$ws.Range("AS2").Value = 200
$ws.Range("AT2").Value = 2.52
$ws.Range("AU2").Value = 6.9
$ws.Range("AV2").Value = 60
$ws.Range("AW2").Value = 5.3
$ws.Range("AX2").Value = 19.5
$ws.Range("AY2").Value = 25
$ws.Range("AZ2").Value = 100
$ws.Range("BA2").Value = 150
$ws.Range("BB2").Value = 300

# Row 3 updates
$ws.Range("G3").Value = 4.1
$ws.Range("H3").Value = 3.65
$ws.Range("I3").Value = 1.7
$ws.Range("J3").Value = 4.5
$ws.Range("L3").Value = 2.25
$ws.Range("Q3").Value = 1.8
$ws.Range("R3").Value = 1.82
$ws.Range("S3").Value = 1.36
$ws.Range("T3").Value = 2.52
$ws.Range("U3").Value = 1.79
$ws.Range("V3").Value = 1.98
$ws.Range("W3").Value = 10
$ws.Range("X3").Value = 18.5
$ws.Range("Y3").Value = 11.75
$ws.Range("Z3").Value = 50
$ws.Range("AA3").Value = 30
$ws.Range("AB3").Value = 35
$ws.Range("AD3").Value = 6.3
$ws.Range("AE3").Value = 13
$ws.Range("AF3").Value = 55
$ws.Range("AG3").Value = 450
$ws.Range("AH3").Value = 6
$ws.Range("AI3").Value = 6.7
$ws.Range("AK3").Value = 10.5
$ws.Range("AL3").Value = 11
$ws.Range("AN3").Value = 5.9
$ws.Range("AO3").Value = 24
$ws.Range("AP3").Value = 32
$ws.Range("AR3").Value = 175
$ws.Range("AS3").Value = 450
$ws.Range("AT3").Value = 2.52
$ws.Range("AU3").Value = 7.7
$ws.Range("AV3").Value = 75
$ws.Range("AW3").Value = 3.5
$ws.Range("AX3").Value = 8.25
$ws.Range("AZ3").Value = 28
